$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N, shifting N:P -> O:Q
$ws.Columns("N").Insert()

# The newly inserted column inherits the width of the column to its left (M)
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet and update its selection/view
$ws.Activate()
$ws.Range("P6").Select() | Out-Null
